$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D (Price) stays text-formatted so numeric-looking
# strings like "67.292.92" or "0.999" are not auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "67.292.92"
$ws.Range("E2").Value = "  +5.84%  "

$ws.Range("D3").Value = "3.715.54"
$ws.Range("E3").Value = "  +7.09%  "

$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").Value = "423.87"
$ws.Range("E5").Value = "  +2.19%  "

$ws.Range("D6").Value = "131.51"
$ws.Range("E6").Value = "  +1.81%  "

$ws.Range("D7").Value = "3.707.06"
$ws.Range("E7").Value = "  +7.09%  "

$ws.Range("E8").Value = "  +2.55%  "

$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  -0.06%  "

$ws.Range("D10").Value = "0.768"
$ws.Range("E10").Value = "  -2.12%  "

$ws.Range("D11").Value = "0.185"
$ws.Range("E11").Value = "  +13.35%  "

$ws.Range("D12").Value = "0.0000398"
$ws.Range("E12").Value = "  +55.30%  "

$ws.Range("D13").Value = "42.99"
$ws.Range("E13").Value = "  +1.31%  "

$ws.Range("D14").Value = "10.15"
$ws.Range("E14").Value = "  +2.99%  "

$ws.Range("D15").Value = "4.298.39"
$ws.Range("E15").Value = "  +6.97%  "

$ws.Range("E16").Value = "  -0.05%  "

$ws.Range("D17").Value = "20.80"
$ws.Range("E17").Value = "  +2.83%  "

$ws.Range("D18").Value = "3.730.17"
$ws.Range("E18").Value = "  +7.04%  "

$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").Value = "12.93"
$ws.Range("E19").Value = "  +4.00%  "

$ws.Range("B20").Value = "Polygon"
$ws.Range("C20").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D20").Value = "1.13"
$ws.Range("E20").Value = "  +4.23%  "

$ws.Range("D21").Value = "67.261.41"
$ws.Range("E21").Value = "  +6.01%  "

$ws.Range("D22").Value = "450.77"
$ws.Range("E22").Value = "  -2.79%  "

$ws.Range("D23").Value = "15.78"
$ws.Range("E23").Value = "  +16.43%  "

$ws.Range("D24").Value = "89.47"
$ws.Range("E24").Value = "  -0.86%  "

$ws.Range("D25").Value = "3.19"
$ws.Range("E25").Value = "  -3.01%  "

$ws.Range("D26").Value = "38.07"
$ws.Range("E26").Value = "  +12.28%  "

$ws.Range("D27").Value = "10.24"
$ws.Range("E27").Value = "  +0.95%  "

$ws.Range("E28").Value = "  +1.50%  "

$ws.Range("D29").Value = "4.98"
$ws.Range("E29").Value = "  +4.59%  "

$ws.Range("D30").Value = "12.70"
$ws.Range("E30").Value = "  +2.21%  "

$ws.Range("E31").Value = "  +9.56%  "

$ws.Range("E32").Value = "  +4.68%  "

$ws.Range("D33").Value = "7.34"
$ws.Range("E33").Value = "  -1.93%  "

$ws.Range("D34").Value = "42.14"
$ws.Range("E34").Value = "  +5.84%  "

$ws.Range("D35").Value = "0.164"
$ws.Range("E35").Value = "  -1.27%  "

$ws.Range("E36").Value = "  -0.01%  "

$ws.Range("D37").Value = "56.45"
$ws.Range("E37").Value = "  -1.93%  "

$ws.Range("D38").Value = "0.0492"
$ws.Range("E38").Value = "  +1.03%  "

$ws.Range("D39").Value = "0.0₃0772"
$ws.Range("E39").Value = "  +17.19%  "

$ws.Range("D40").Value = "3.17"
$ws.Range("E40").Value = "  +36.22%  "

$ws.Range("E41").Value = "  +5.26%  "

$ws.Range("D42").Value = "28.24"
$ws.Range("E42").Value = "  +29.55%  "

$ws.Range("D43").Value = "0.998"
$ws.Range("E43").Value = "  -0.06%  "

$ws.Range("D44").Value = "3.43"
$ws.Range("E44").Value = "  +3.21%  "

$ws.Range("D45").Value = "2.95"
$ws.Range("E45").Value = "  -4.06%  "

$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").Value = "4.44"
$ws.Range("E46").Value = "  -1.45%  "

$ws.Range("B47").Value = "Monero"
$ws.Range("C47").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D47").Value = "146.53"
$ws.Range("E47").Value = "  +1.63%  "

$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").Value = "2.11"
$ws.Range("E48").Value = "  +5.65%  "

$ws.Range("D49").Value = "2.69"
$ws.Range("E49").Value = "  -2.93%  "

$ws.Range("D50").Value = "0.309"
$ws.Range("E50").Value = "  -2.30%  "

$ws.Range("D51").Value = "0.159"
$ws.Range("E51").Value = "  +16.56%  "
